$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196. This shifts the existing rows 196-220
# down to 197-221 (preserving all their data/formatting), matching the
# target diff which adds one new record and pushes the rest of the table
# down by one row.
$ws.Rows(196).Insert()

# Populate the newly inserted row 196 with the new weekly price record.
$ws.Cells.Item(196, 1).Value  = 9
$ws.Cells.Item(196, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(196, 3).Value  = "Metropolitana"
$ws.Cells.Item(196, 4).Value  = 44505
$ws.Cells.Item(196, 5).Value  = 13
$ws.Cells.Item(196, 6).Value  = 100112052
$ws.Cells.Item(196, 7).Value  = "Albahaca"
$ws.Cells.Item(196, 8).Value  = "Sin especificar"
$ws.Cells.Item(196, 9).Value  = "Primera"
$ws.Cells.Item(196, 10).Value = 70
$ws.Cells.Item(196, 11).Value = 6000
$ws.Cells.Item(196, 12).Value = 6000
$ws.Cells.Item(196, 13).Value = 6000
$ws.Cells.Item(196, 14).Value = "`$/docena de matas"
$ws.Cells.Item(196, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(196, 16).Value = 1000
$ws.Cells.Item(196, 17).Value = 6
$ws.Cells.Item(196, 18).Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date/time number format,
# same as the rest of column D.
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
